$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "62.766.27"
$ws.Cells.Item(2, 5).Value = "  +3.27%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.447.50"
$ws.Cells.Item(3, 5).Value = "  +2.03%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.12%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "577.98"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.92%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "145.63"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +2.82%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.03%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.14%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "2.445.94"
$ws.Cells.Item(9, 5).Value = "  +1.74%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +3.11%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +1.94%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +1.15%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.354"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +2.57%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +7.74%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.0000179"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +5.60%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.890.24"

# Row 17
$ws.Cells.Item(17, 4).Value = "62.738.22"
$ws.Cells.Item(17, 5).Value = "  +3.84%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "2.452.72"
$ws.Cells.Item(18, 5).Value = "  +1.43%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.85"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -3.35%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "10.96"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +2.80%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "0.0₆0843"
$ws.Cells.Item(21, 5).Value = "  +201.56%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "330.33"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +1.90%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "4.14"
$ws.Cells.Item(23, 4).Style = "Normal"

# Row 24
$ws.Cells.Item(24, 5).Value = "  +9.19%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +0.14%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "65.91"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.97%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "645.60"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +12.67%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +17.17%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "8.46"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +4.42%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "0.0₃0985"
$ws.Cells.Item(30, 5).Value = "  +4.79%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +9.01%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +1.59%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +3.58%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +4.27%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +2.58%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.999"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.11%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "5.53"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +7.17%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Monero"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "153.21"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.69%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.374"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.91%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "18.73"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +2.33%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.71"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +6.40%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +4.69%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "USDe"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.999"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.02%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "OKB"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "42.21"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.37%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "14.93"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +27.37%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "145.33"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +2.24%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "3.61"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +2.47%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "20.65"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +6.64%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  +3.03%  "
